$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Morocco ✓ - Republic of the Congo: 1:0"
$ws.Range("B2").Value = "Morocco"
$ws.Range("C2").Value = 82
$ws.Range("D2").Value = 80
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = "✓"

# Row 3
$ws.Range("A3").Value = "Latvia - England ✓: 0:5"
$ws.Range("B3").Value = "England"
$ws.Range("C3").Value = 80
$ws.Range("D3").Value = 97
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = 1.05
$ws.Range("G3").Value = "✓"

# Row 4
$ws.Range("A4").Value = "Spain ✓ - Bulgaria: 4:0"
$ws.Range("B4").Value = "Spain"
$ws.Range("C4").Value = 79
$ws.Range("D4").Value = 86
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 1.02
$ws.Range("G4").Value = "✓"

# Row 5 (new)
$ws.Range("A5").Value = "Portugal  - Hungary: 2:2"
$ws.Range("B5").Value = "Portugal"
$ws.Range("C5").Value = 77
$ws.Range("D5").Value = 95
$ws.Range("E5").Value = 97
$ws.Range("F5").Value = 1.23
$ws.Range("G5").Value = ""

# Row 6 (new)
$ws.Range("A6").Value = "Senegal ✓ - Mauritania: 4:0"
$ws.Range("B6").Value = "Senegal"
$ws.Range("C6").Value = 74
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 1.15
$ws.Range("G6").Value = "✓"

# Row 7 (new)
$ws.Range("A7").Value = "Italy ✓ - Israel: 3:0"
$ws.Range("B7").Value = "Italy"
$ws.Range("C7").Value = 74
$ws.Range("D7").Value = 87
$ws.Range("E7").Value = 97
$ws.Range("F7").Value = 1.22
$ws.Range("G7").Value = "✓"

# Row 8 (new) - was old row 3, shifted down with updates
$ws.Range("A8").Value = "Democratic Republic of the Congo ✓ - Sudan: 1:0"
$ws.Range("B8").Value = "Democratic Republic of the Congo"
$ws.Range("C8").Value = 73
$ws.Range("D8").Value = 79
$ws.Range("E8").Value = 73
$ws.Range("F8").Value = 1.57
$ws.Range("G8").Value = "✓"

# Row 9 (new) - was old row 4, shifted down with updates
$ws.Range("A9").Value = "Ivory Coast ✓ - Kenya: 3:0"
$ws.Range("B9").Value = "Ivory Coast"
$ws.Range("C9").Value = 73
$ws.Range("D9").Value = 77
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = 1.11
$ws.Range("G9").Value = "✓"

# Row 10 (new)
$ws.Range("A10").Value = "Gabon ✓ - Burundi: 2:0"
$ws.Range("B10").Value = "Gabon"
$ws.Range("C10").Value = 73
$ws.Range("D10").Value = 89
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = 1.3
$ws.Range("G10").Value = "✓"

# Row 11 (new)
$ws.Range("A11").Value = "Turkiye ✓ - Georgia: 4:1"
$ws.Range("B11").Value = "Turkiye"
$ws.Range("C11").Value = 57
$ws.Range("D11").Value = 85
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 1.57
$ws.Range("G11").Value = "✓"

"Done"
